# NCATS Study multifilter testcases61to70
# - Update the Cypher query stored in B4 (StudyFilesTab query): change the
#   WHERE clause filter from demo.sex/file_type/tif to
#   diag.stage_of_disease/rtf (file_type filter removed).
# - Move the active sheet selection from C4 to A5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$origHeight = $ws.Rows.Item(4).RowHeight

$ws.Range("B4").Value = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nMATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nOPTIONAL MATCH (f)-->(samp:sample)`nMATCH (f)-->(diag:diagnosis)`nWHERE s.clinical_study_designation IN ['NCATS-COP01'] and diag.stage_of_disease in ['III'] and labels(parent)[0] IN ['diagnosis'] and f.file_format IN ['rtf']`nWITH`n        DISTINCT f, parent, c, demo, diag, s, samp,`n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`n        2 as precision`nWITH`n        f, parent, c, demo, diag, s, samp,`n        f.file_size /(1024^i) AS value,`n        10^precision AS factor,`n        units[i] as unit`nWITH`n        f, parent, c, demo, diag, s, samp, unit,`n        round(factor * value)/factor AS size`nRETURN`n        coalesce(f.file_name, '') AS ``File Name``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_type, '') AS ``File Type``,`n        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(samp.sample_id, '') AS ``Sample ID``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n        coalesce(demo.breed,'') AS Breed ,`n        coalesce(diag.disease_term,'') AS Diagnosis`n        order by f.file_name asc`n        limit 100"

# Setting a much longer value triggers an automatic row autofit in this
# engine; restore the original explicit row height so it stays unchanged.
$ws.Rows.Item(4).RowHeight = $origHeight

$ws.Range("A5").Select()
